$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.407.36"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "2.011.75"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'259.82"
$ws.Range("E5").Value = "  +4.69%  "
$ws.Range("D6").Value = "'0.620"
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'57.17"
$ws.Range("E8").Value = "  -5.69%  "
$ws.Range("D9").Value = "'0.384"
$ws.Range("E9").Value = "  -3.61%  "
$ws.Range("D10").Value = "'0.0776"
$ws.Range("E10").Value = "  -4.47%  "
$ws.Range("E11").Value = "  -3.16%  "
$ws.Range("D12").Value = "'14.31"
$ws.Range("D13").Value = "2.305.24"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").Value = "'21.53"
$ws.Range("E14").Value = "  -4.24%  "
$ws.Range("D15").Value = "'0.800"
$ws.Range("E15").Value = "  -7.56%  "
$ws.Range("D16").Value = "'5.24"
$ws.Range("E16").Value = "  -5.35%  "
$ws.Range("D17").Value = "2.031.43"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "37.247.78"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "'69.89"
$ws.Range("E19").Value = "  -1.46%  "
$ws.Range("D20").Value = "0.0₃0839"
$ws.Range("E20").Value = "  -3.54%  "
$ws.Range("D21").Value = "'232.09"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "'5.13"
$ws.Range("E22").Value = "  -2.88%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "'2.59"
$ws.Range("E23").Value = "  +3.24%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").Value = "'165.01"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").Value = "'8.97"
$ws.Range("E27").Value = "  -5.61%  "
$ws.Range("D28").Value = "'19.60"
$ws.Range("E28").Value = "  -1.36%  "
$ws.Range("D29").Value = "'0.129"
$ws.Range("E29").Value = "  -6.51%  "
$ws.Range("E30").Value = "  -5.12%  "
$ws.Range("D31").Value = "'0.120"
$ws.Range("E31").Value = "  -2.20%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.60"
$ws.Range("E32").Value = "  -5.40%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0643"
$ws.Range("E33").Value = "  -3.58%  "
$ws.Range("D34").Value = "'4.51"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("D35").Value = "'2.37"
$ws.Range("E35").Value = "  -5.82%  "
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "'3.36"
$ws.Range("E38").Value = "  -2.84%  "
$ws.Range("D39").Value = "'5.48"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("D40").Value = "'3.05"
$ws.Range("E40").Value = "  +3.34%  "
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").Value = "'0.0931"
$ws.Range("E42").Value = "  -5.25%  "
$ws.Range("D43").Value = "'0.0212"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("D44").Value = "1.420.89"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("D45").Value = "'15.79"
$ws.Range("E45").Value = "  -7.11%  "
$ws.Range("D46").Value = "'89.69"
$ws.Range("E47").Value = "  -3.76%  "
$ws.Range("D48").Value = "'2.93"
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("D49").Value = "'7.02"
$ws.Range("E49").Value = "  -7.44%  "
$ws.Range("D50").Value = "2.196.47"
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("E51").Value = "  -9.65%  "
